$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1. Stage")

$ws.Range("I2").Value = 0.073742
$ws.Range("I5").Value = 0.362848
$ws.Range("I6").Value = 0.220403
$ws.Range("I7").Value = 0.757201
$ws.Range("I10").Value = 0.000894
$ws.Range("I11").Value = 0.000125
$ws.Range("I12").Value = 0.001877
$ws.Range("I13").Formula = "=SUM(I2:I12)"

$ws.Range("I13").Select()
